$d = $word.ActiveDocument

# Insert a new paragraph right after the title ("ETL Project Report")
# and before "Extract:" to hold the authors' names.
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$namesPara = $d.Paragraphs(2)

# Make it bold (both the paragraph mark run props and the text run props)
# before typing the text so the new run inherits the bold formatting.
$namesPara.Range.Font.Bold = 1
$namesPara.Range.Font.BoldBi = 1
$namesPara.Range.Text = "By Achyut Shrestha, Kiran Babuk, Sam Deschamps, Mark Gotanga"

# The "_GoBack" bookmark (Word's "last edit location" marker) moves from
# the end of the document to this newly-added paragraph.
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

$namesPara2 = $d.Paragraphs(2)
$paraStart = $namesPara2.Range.Start
$paraEnd = $namesPara2.Range.End
# Exclude the trailing paragraph mark so the bookmark wraps just the text.
$bookmarkRange = $d.Range($paraStart, $paraEnd - 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
